# Tutorial 6 solution update:
#  - Reformat the Date column (A3:A21) from DD/MM/YYYY to DD-MM-YYYY
#  - Correct a few attendance tally cells that depended on the date parsing

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New dash-separated dates for rows 3-21 (column A)
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

# Excel auto-parses dash-separated day/month/year text as a date when the
# first number is <= 12 (ambiguous as US-style mm-dd-yyyy). Force the cell
# to Text before writing so the value is kept as a literal string, then
# clear the formatting override so no stray style index is left behind
# (matches the source file, where these cells carry no explicit style).
foreach ($row in ($dates.Keys | Sort-Object)) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$row]
    $cell.ClearFormats()
}

# Attendance tally corrections
# Row 3: Total Attendance Count (D) and Invalid (G) go from 0 to 1
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 7).Value = 1

# Row 6: Total Attendance Count (D) and Real (E) go from 0 to 1, Absent (H) goes from 1 to 0
$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 8).Value = 0
